# Updated the storage class diagram of developer guide
#
# The slide's storage/"Budget" elbow connectors were drawn with a
# purple (7030A0) system-dot dashed outline. Restyle them to a solid
# single line (matches PowerPoint's "Style" -> single compound line,
# which writes cmpd="sng", and "DashStyle" -> solid, which rewrites
# prstDash from "sysDot" to "solid").
#
# Shapes.Item(62) -> cxnSp id="108" "Elbow Connector 63" (adj1 = 50000)
# Shapes.Item(64) -> cxnSp id="110" "Elbow Connector 63" (adj1 = -124816)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$connectors = @(62, 64)

foreach ($idx in $connectors) {
    $shape = $s.Shapes.Item($idx)
    $line = $shape.Line

    # msoLineSingle -> adds cmpd="sng" to <a:ln>
    $line.Style = 1
    # msoLineSolid -> <a:prstDash val="solid"/> (was "sysDot")
    $line.DashStyle = 1
}
